# Update the cryptocurrency list on the active worksheet with the latest
# price/volume(1h) figures from the scheduled GitHub Actions refresh.
# Includes two coin-row swaps: ImmutableX/VeChain (rows 37-38) and
# Kaspa/FraxShare (rows 44-45). Column D values are re-applied as text
# (via a temporary "@" number format) so that figures such as "226.99"
# or "0.0509" remain plain text rather than being parsed as numbers,
# matching how the sheet stores them as inline strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '34.852.55'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +3.00%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.807.08'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.34%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.997'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.59%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '226.99'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.92%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.559'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.63%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.996'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.44%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '33.33'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +6.94%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.288'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +3.02%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0675'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.57%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0937'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.39%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '2.063.23'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.23%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '11.57'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +16.32%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.792.71'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.50%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.646'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.89%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '34.776.59'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +3.03%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '4.32'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +3.42%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '69.90'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.80%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '258.12'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.60%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0773'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +4.85%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.998'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.71%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.55'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.53%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.28'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.36%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.14'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.43%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '159.06'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.63%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '16.61'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.25%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.20'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +4.16%  '
$ws.Range('E28').Value = '  +0.32%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.997'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.42%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.84'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.97%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0525'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +3.11%  '
$ws.Range('E32').Value = '  +0.89%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.63'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.16%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.93'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +11.02%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.468.63'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.29%  '
$ws.Range('E36').Value = '  +0.73%  '
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.641'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +3.25%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0192'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.38%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '84.23'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.25%  '
$ws.Range('E40').Value = '  +3.71%  '
$ws.Range('E41').Value = '  -1.19%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.911'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.22%  '
$ws.Range('E43').Value = '  +1.79%  '
$ws.Range('B44').Value = 'Kaspa'
$ws.Range('C44').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0509'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.04%  '
$ws.Range('B45').Value = 'FraxShare'
$ws.Range('C45').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '6.04'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +6.21%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.962.59'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.67%  '
$ws.Range('E47').Value = '  -2.74%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '12.12'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.02%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.998'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.68%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '103.31'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +4.85%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '50.66'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.17%  '
